$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 (date 42880 / 2017-05-25) ---
$ws.Range("B18:J18").Value = "X"
$ws.Range("K18").Value = "Termino da Introdução da parte escrita Tcc"
$ws.Range("L18").Value = "Atualizado o Kanban"
$ws.Rows.Item(18).RowHeight = 16.5

# --- Row 19 (date 42887 / 2017-06-01) ---
$ws.Range("A19").Value = 42887
$ws.Range("B19:J19").Value = "X"
$ws.Range("K19").Value = " continuação da Aula sobre formatação ABNT no word"
$ws.Rows.Item(19).RowHeight = 16.5

# --- Row 20 (date 42894 / 2017-06-08) ---
$ws.Range("A20").Value = 42894
$ws.Range("B20:J20").Value = "X"
$ws.Range("K20").Value = "Reuniao de ajustes do mockup, tirado duvidas sobre a estrutura do site"
$ws.Rows.Item(20).RowHeight = 30

# --- Update the active selection to K20, matching the saved view state ---
$ws.Range("K20").Select()
